# Banshee_UDP_Mapping.xlsx edits
# 1) Fixed F4 Ess&PV relay settings
# 2) Modified sync check to reuse pll measurements already available in the model
# 3) Updated udp frame (unified frame for PHIL and Banshee)
#
# The "PacketFormat" sheet lists the UDP frame layout; column E holds the
# per-point byte count ("Qty") and column F a shared formula (=D*E) that
# derives the total bytes for that row. Rows 6:10 (the F4 Ess/PV relay
# block) move from 48 -> 59 bytes/point, and rows 35:36 (the sync-check /
# PLL measurement rows) move from 18 -> 25 bytes/point. The running totals
# in E58 and F60 are formulas and recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PacketFormat")

# 1) F4 Ess & PV relay settings - Qty column (E) for rows 6-10
$ws.Range("E6:E10").Value = 59

# 3) Updated UDP frame - sync check / PLL measurement rows 35-36
$ws.Range("E35:E36").Value = 25

# Keep the active sheet/selection pointed at the area that was edited,
# matching where the author ended up after making the change.
$ws.Activate()
$ws.Range("E37").Select()
